# Add new weekly scoreboard rows (111-115) to Sheet1, matching the
# "Kilimanjaro Weekly Scoreboard" workbook update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the date number format (style) from the last existing date cell
# (B110) onto the new date cells, so they render/serialize the same way
# ("m/d/yyyy"-style numFmtId=14) as the rest of column B.
$ws.Range("B110").Copy()
$ws.Range("B111:B115").PasteSpecial(-4122)

# --- Row 111: Steven, Workout, Wily Hyena
$ws.Cells.Item(111, 1).Value = "Steven"
$ws.Cells.Item(111, 2).Value = 45468
$ws.Cells.Item(111, 3).Value = "Workout"
$ws.Cells.Item(111, 4).Value = 33
$ws.Cells.Item(111, 5).Value = 0
$ws.Cells.Item(111, 6).Value = 0
$ws.Cells.Item(111, 7).Value = 21
$ws.Cells.Item(111, 8).Value = 11
$ws.Cells.Item(111, 9).Value = 1
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 11).Value = 0
$ws.Cells.Item(111, 12).Value = "Wily Hyena"
$ws.Cells.Item(111, 13).Value = 3

# --- Row 112: Eric, Workout, Sauntering Hippo (new team name)
$ws.Cells.Item(112, 1).Value = "Eric"
$ws.Cells.Item(112, 2).Value = 45468
$ws.Cells.Item(112, 3).Value = "Workout"
$ws.Cells.Item(112, 4).Value = 73
$ws.Cells.Item(112, 5).Value = 0
$ws.Cells.Item(112, 6).Value = 0
$ws.Cells.Item(112, 7).Value = 8
$ws.Cells.Item(112, 8).Value = 52
$ws.Cells.Item(112, 9).Value = 12
$ws.Cells.Item(112, 10).Value = 3
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 12).Value = "Sauntering Hippo"
$ws.Cells.Item(112, 13).Value = 3

# --- Row 113: Steven, Walk, Wily Hyena
$ws.Cells.Item(113, 1).Value = "Steven"
$ws.Cells.Item(113, 2).Value = 45468
$ws.Cells.Item(113, 3).Value = "Walk"
$ws.Cells.Item(113, 4).Value = 35
$ws.Cells.Item(113, 5).Value = 1.63
$ws.Cells.Item(113, 6).Value = 43
$ws.Cells.Item(113, 7).Value = 35
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = "Wily Hyena"
$ws.Cells.Item(113, 13).Value = 3

# --- Row 114: Matt, Ride, Agile Antelope
$ws.Cells.Item(114, 1).Value = "Matt"
$ws.Cells.Item(114, 2).Value = 45469
$ws.Cells.Item(114, 3).Value = "Ride"
$ws.Cells.Item(114, 4).Value = 79
$ws.Cells.Item(114, 5).Value = 21.21
$ws.Cells.Item(114, 6).Value = 1184
$ws.Cells.Item(114, 7).Value = 6
$ws.Cells.Item(114, 8).Value = 67
$ws.Cells.Item(114, 9).Value = 5
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 11).Value = 0
$ws.Cells.Item(114, 12).Value = "Agile Antelope"
$ws.Cells.Item(114, 13).Value = 3

# --- Row 115: Jeremiah, Run, Agile Antelope
$ws.Cells.Item(115, 1).Value = "Jeremiah"
$ws.Cells.Item(115, 2).Value = 45469
$ws.Cells.Item(115, 3).Value = "Run"
$ws.Cells.Item(115, 4).Value = 19
$ws.Cells.Item(115, 5).Value = 2.14
$ws.Cells.Item(115, 6).Value = 157
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 3
$ws.Cells.Item(115, 9).Value = 10
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 11).Value = 0
$ws.Cells.Item(115, 12).Value = "Agile Antelope"
$ws.Cells.Item(115, 13).Value = 3

# --- Column L ("Workout Level") needs to re-fit its width now that the
# longer "Sauntering Hippo" label is present (bestFit grows 13.55 -> 14.78
# character-width units in the original; closest value reachable through
# the COM ColumnWidth API, which only resolves to 1/6-character
# increments, is used here).
$ws.Columns.Item(12).ColumnWidth = 13.95

# --- Move the frozen-pane viewport / active selection down to the new
# bottom of the data, as Excel does after appending rows and navigating
# to the new last row.
$ws.Range("A116").Select()
